# Updated cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.431.48"
$ws.Range("E2").Value = "  -1.46%  "
$ws.Range("D3").Value = "1.842.39"
$ws.Range("E3").Value = "  -1.78%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "261.35"
$ws.Range("E5").Value = "  -5.72%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5212"
$ws.Range("E7").Value = "  -0.97%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3269"
$ws.Range("E8").Value = "  -4.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06783"
$ws.Range("E9").Value = "  -2.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.62"
$ws.Range("E10").Value = "  -7.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7714"
$ws.Range("E11").Value = "  -3.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07696"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").Value = "1.837.20"
$ws.Range("E13").Value = "  -2.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.09"
$ws.Range("E14").Value = "  -2.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.008"
$ws.Range("E15").Value = "  -3.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.92"
$ws.Range("E17").Value = "  -4.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007973"
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("D20").Value = "26.438.51"
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("D21").Value = "2.067.02"
$ws.Range("E21").Value = "  -1.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.580"
$ws.Range("E22").Value = "  -3.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.484"
$ws.Range("E23").Value = "  -5.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.978"
$ws.Range("E24").Value = "  -3.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.66"
$ws.Range("E25").Value = "  -1.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.198"
$ws.Range("E26").Value = "  -7.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.650"
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("E28").Value = "  -2.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.52"
$ws.Range("E29").Value = "  -1.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.203"
$ws.Range("E30").Value = "  -3.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.128"
$ws.Range("E31").Value = "  -4.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08707"
$ws.Range("E32").Value = "  -2.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04790"
$ws.Range("E33").Value = "  -2.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.131"
$ws.Range("E34").Value = "  -3.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7162"
$ws.Range("E35").Value = "  -1.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.842"
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.081"
$ws.Range("E37").Value = "  -6.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.223"
$ws.Range("E38").Value = "  -4.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01771"
$ws.Range("E39").Value = "  -3.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4838"
$ws.Range("E40").Value = "  -5.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "111.97"
$ws.Range("E41").Value = "  -3.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8978"
$ws.Range("E42").Value = "  -6.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.064"
$ws.Range("E43").Value = "  -1.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9999"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.719"
$ws.Range("E45").Value = "  -4.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4152"
$ws.Range("E46").Value = "  -6.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05877"
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.966"
$ws.Range("E48").Value = "  -4.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.01"
$ws.Range("E49").Value = "  -3.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1216"
$ws.Range("E50").Value = "  -9.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8885"
$ws.Range("E51").Value = "  +0.37%  "
